$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.509.44"

$ws.Range("D3").Value = "1.675.90"
$ws.Range("E3").Value = "  +1.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.99"
$ws.Range("E5").Value = "  +1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5318"
$ws.Range("E6").Value = "  +2.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2696"
$ws.Range("E8").Value = "  +3.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06404"
$ws.Range("E9").Value = "  +0.37%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.80"
$ws.Range("E10").Value = "  +4.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +1.45%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.693.17"
$ws.Range("E12").Value = "  +3.31%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.508"
$ws.Range("E13").Value = "  +1.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5582"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").Value = "0.0₅8362"
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.76"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").Value = "26.525.31"
$ws.Range("E17").Value = "  +1.64%  "

$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.779"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.01"
$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.35"
$ws.Range("E21").Value = "  +1.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.326"
$ws.Range("E22").Value = "  +1.53%  "

$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1284"
$ws.Range("E24").Value = "  +5.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.63"
$ws.Range("E25").Value = "  -4.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.432"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.30"
$ws.Range("E27").Value = "  +3.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.431"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06291"
$ws.Range("E29").Value = "  +6.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.294"
$ws.Range("E30").Value = "  +1.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.605"
$ws.Range("E31").Value = "  +5.61%  "

$ws.Range("E32").Value = "  +1.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.696"
$ws.Range("E33").Value = "  +2.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.014"
$ws.Range("E34").Value = "  +2.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6169"
$ws.Range("E35").Value = "  +9.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.425"
$ws.Range("E36").Value = "  +1.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.786"
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.122"
$ws.Range("E39").Value = "  +4.97%  "

$ws.Range("D40").Value = "1.095.55"
$ws.Range("E40").Value = "  +6.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8613"
$ws.Range("E41").Value = "  +0.35%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.65"
$ws.Range("E43").Value = "  +0.48%  "

$ws.Range("D44").Value = "1.822.92"
$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("E45").Value = "  +4.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.84"
$ws.Range("E46").Value = "  +5.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.202"
$ws.Range("E47").Value = "  +1.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  +0.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.514"
$ws.Range("E49").Value = "  +9.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05194"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.028"
$ws.Range("E51").Value = "  +1.65%  "
